# Text updates as supplied by PM&C.
# Adds "Source" and "References" information to the bottom of the
# Description sheet (rows 10-12).

$wb = $excel.ActiveWorkbook

# Reset the (inactive) Data sheet's lingering selection back to A1.
$dataWs = $wb.Worksheets.Item("Data")
$dataWs.Range("A1").Select() | Out-Null

$ws = $wb.Worksheets.Item("Description")
$ws.Activate() | Out-Null

# Row 10: Source
$ws.Range("A10").Value = "Source"
$ws.Range("B10").Value = "ABS (various years) Education and Work, Australia."

# Row 11: References (heading + first reference)
$ws.Range("A11").Value = "References"
$ws.Range("B11").Value = "Beddie, F. (2015). The outcomes of education and training: what the Australian research is telling us, 2011-14. Adelaide: NCVER."

# Row 12: second reference (column B only)
$ws.Range("B12").Value = "Wheelahan, L., Buchanan, J. and Yu, S. (2015). Linking qualifications and the labour market through capabilities and vocational streams. Adelaide: NCVER."

# Match the wrapped-text look used by the rest of column B in this sheet.
$ws.Range("B10:B12").WrapText = $true

# Row heights sized to fit the wrapped reference text.
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 26.85
$ws.Rows.Item(12).RowHeight = 26.95

$ws.Range("B12").Select() | Out-Null
